# Auto-generated edit script applying scheduled market-data refresh
# to the Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 108.07692
$ws.Range("I9").Value = 105.454544
$ws.Range("J9").Value = 122.5
$ws.Range("K9").Value = 105.454544
$ws.Range("L9").Value = 122.5
$ws.Range("M9").Value = 63.545456
$ws.Range("N9").Value = -460.5
$ws.Range("H11").Value = 149.27272
$ws.Range("I11").Value = 149.27272
$ws.Range("K11").Value = 149.27272
$ws.Range("M11").Value = -9.272719999999993
$ws.Range("H12").Value = 355.75
$ws.Range("J12").Value = 206.5
$ws.Range("L12").Value = 206.5
$ws.Range("N12").Value = -546.5
$ws.Range("H40").Value = 3485.3333
$ws.Range("J40").Value = 3485.3333
$ws.Range("L40").Value = 3485.3333
$ws.Range("N40").Value = -3835.3333
$ws.Range("H107").Value = 2166.889
$ws.Range("I107").Value = 1037.9231
$ws.Range("J107").Value = 5102.2
$ws.Range("K107").Value = 1037.9231
$ws.Range("L107").Value = 5102.2
$ws.Range("M107").Value = 882.0769
$ws.Range("N107").Value = -8942.200000000001
$ws.Range("H137").Value = 2574.2
$ws.Range("I137").Value = 1843.3334
$ws.Range("J137").Value = 3670.5
$ws.Range("K137").Value = 5530.0002
$ws.Range("L137").Value = 11011.5
$ws.Range("M137").Value = -2980.0002
$ws.Range("N137").Value = -16111.5
$ws.Range("H138").Value = 1111.5454
$ws.Range("I138").Value = 837.7
$ws.Range("K138").Value = 2513.1
$ws.Range("M138").Value = 2626.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4556.115
$ws.Range("I32").Value = 3665.0588
$ws.Range("K32").Value = 3665.0588
$ws.Range("M32").Value = -3378.0588
$ws.Range("H36").Value = 6000
$ws.Range("I36").Value = 5000
$ws.Range("J36").Value = 7000
$ws.Range("K36").Value = 5000
$ws.Range("L36").Value = 7000
$ws.Range("M36").Value = -4654
$ws.Range("N36").Value = -7692
$ws.Range("H61").Value = 1383.5714
$ws.Range("I61").Value = 950
$ws.Range("K61").Value = 950
$ws.Range("M61").Value = -738
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("H132").Value = 8616.375
$ws.Range("I132").Value = 8986.166999999999
$ws.Range("K132").Value = 26958.501
$ws.Range("M132").Value = -24428.501
$ws.Range("H136").Value = 1383.5714
$ws.Range("I136").Value = 950
$ws.Range("K136").Value = 2850
$ws.Range("M136").Value = -300
$ws.Range("N62").ClearContents()
$ws.Range("N65").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1825.8235
$ws.Range("I20").Value = 1656.75
$ws.Range("K20").Value = 1656.75
$ws.Range("M20").Value = -1409.75
$ws.Range("H22").Value = 125.14286
$ws.Range("I22").Value = 79.333336
$ws.Range("J22").Value = 400
$ws.Range("K22").Value = 79.333336
$ws.Range("L22").Value = 400
$ws.Range("M22").Value = 93.666664
$ws.Range("N22").Value = -746
$ws.Range("H107").Value = 3173

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 705.1667
$ws.Range("I19").Value = 236.57143
$ws.Range("J19").Value = 1361.2
$ws.Range("K19").Value = 236.57143
$ws.Range("L19").Value = 1361.2
$ws.Range("M19").Value = -66.57142999999999
$ws.Range("N19").Value = -1701.2
$ws.Range("H22").Value = 457.7143
$ws.Range("I22").Value = 112.5
$ws.Range("K22").Value = 112.5
$ws.Range("M22").Value = 237.5
$ws.Range("H24").Value = 705.1667
$ws.Range("I24").Value = 236.57143
$ws.Range("J24").Value = 1361.2
$ws.Range("K24").Value = 236.57143
$ws.Range("L24").Value = 1361.2
$ws.Range("M24").Value = -66.57142999999999
$ws.Range("N24").Value = -1701.2
$ws.Range("H31").Value = 4256
$ws.Range("I31").Value = 3141.5715
$ws.Range("K31").Value = 3141.5715
$ws.Range("M31").Value = -2846.5715
$ws.Range("H34").Value = 4256
$ws.Range("I34").Value = 3141.5715
$ws.Range("K34").Value = 3141.5715
$ws.Range("M34").Value = -2939.5715
$ws.Range("H59").Value = 30039.062
$ws.Range("J59").Value = 34997.5
$ws.Range("L59").Value = 34997.5
$ws.Range("N59").Value = -37287.5
$ws.Range("H60").Value = 24997.223
$ws.Range("J60").Value = 24997.223
$ws.Range("L60").Value = 24997.223
$ws.Range("N60").Value = -26019.223
$ws.Range("H99").Value = 6714.143
$ws.Range("I99").Value = 6999.8335
$ws.Range("K99").Value = 6999.8335
$ws.Range("M99").Value = -5501.8335
$ws.Range("H109").Value = 78306.75
$ws.Range("J109").Value = 86656
$ws.Range("L109").Value = 86656
$ws.Range("N109").Value = -88736
$ws.Range("H126").Value = 6714.143
$ws.Range("I126").Value = 6999.8335
$ws.Range("K126").Value = 20999.5005
$ws.Range("M126").Value = -18529.5005
$ws.Range("H132").Value = 2812.5
$ws.Range("I132").Value = 2812.5
$ws.Range("K132").Value = 8437.5
$ws.Range("M132").Value = -5907.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 8503.5
$ws.Range("J81").Value = 9004
$ws.Range("L81").Value = 27012
$ws.Range("N81").Value = -29258
$ws.Range("H84").Value = 8503.5
$ws.Range("J84").Value = 9004
$ws.Range("L84").Value = 81036
$ws.Range("N84").Value = -92268
$ws.Range("H87").Value = 328
$ws.Range("I87").Value = 328
$ws.Range("K87").Value = 984
$ws.Range("M87").Value = 264
$ws.Range("H90").Value = 328
$ws.Range("I90").Value = 328
$ws.Range("K90").Value = 2952
$ws.Range("M90").Value = 3288
$ws.Range("H120").Value = 19240.715
$ws.Range("I120").Value = 3092.5
$ws.Range("K120").Value = 9277.5
$ws.Range("M120").Value = -4439.5
$ws.Range("H122").Value = 1239.5555
$ws.Range("I122").Value = 1108.4445
$ws.Range("J122").Value = 1370.6666
$ws.Range("K122").Value = 9976.0005
$ws.Range("L122").Value = 12335.9994
$ws.Range("M122").Value = -7526.0005
$ws.Range("N122").Value = -17235.9994
$ws.Range("H132").Value = 1580.5385
$ws.Range("I132").Value = 1006
$ws.Range("K132").Value = 9054
$ws.Range("M132").Value = -6524

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 1195
$ws.Range("J25").Value = 1195
$ws.Range("L25").Value = 1195
$ws.Range("N25").Value = -2253
$ws.Range("I47").Value = 31000
$ws.Range("K47").Value = 31000
$ws.Range("M47").Value = -30432
$ws.Range("H48").Value = 27499.5
$ws.Range("I48").Value = 27499.5
$ws.Range("K48").Value = 27499.5
$ws.Range("M48").Value = -27014.5
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 1750
$ws.Range("I38").Value = 1500
$ws.Range("K38").Value = 1500
$ws.Range("M38").Value = -1090
$ws.Range("H61").Value = 7988.5557
$ws.Range("I61").Value = 6986
$ws.Range("K61").Value = 6986
$ws.Range("M61").Value = -6784
$ws.Range("H99").Value = 43500
$ws.Range("I99").Value = 43500
$ws.Range("K99").Value = 43500
$ws.Range("M99").Value = -40505
$ws.Range("H113").Value = 7988.5557
$ws.Range("I113").Value = 6986
$ws.Range("K113").Value = 6986
$ws.Range("M113").Value = -4816
$ws.Range("H122").Value = 5201.2
$ws.Range("I122").Value = 3335.3333
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 10005.9999
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -7555.999899999999
$ws.Range("N122").Value = -28900
$ws.Range("H136").Value = 3612
$ws.Range("I136").Value = 3425.111
$ws.Range("J136").Value = 3780.2
$ws.Range("K136").Value = 10275.333
$ws.Range("L136").Value = 11340.6
$ws.Range("M136").Value = -7725.332999999999
$ws.Range("N136").Value = -16440.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1376.625
$ws.Range("I113").Value = 861.6
$ws.Range("J113").Value = 2235
$ws.Range("K113").Value = 2584.8
$ws.Range("L113").Value = 6705
$ws.Range("M113").Value = -414.8000000000002
$ws.Range("N113").Value = -11045
$ws.Range("H136").Value = 2906.7646
$ws.Range("I136").Value = 2791.2144
$ws.Range("K136").Value = 8373.643199999999
$ws.Range("M136").Value = -5823.643199999999
